$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The old scraper only pulled team/player statistics; now also capture the
# team's season record (Wins / Losses / Ties) alongside each player row.

# New header cells, matching the existing header formatting (bold, bordered,
# centered) used by the rest of row 1.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Populate every player row (2-56) with the team's season record.
for ($r = 2; $r -le 56; $r++) {
    $ws.Cells.Item($r, 30).Value = 73   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 89   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
